$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.874.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.372.99'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.29'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.15%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.368.53'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.21%  '
$ws.Range('E10').Value = '  -0.11%  '
$ws.Range('E11').Value = '  +1.34%  '
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.960.97'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.96%  '
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('E15').Value = '  -3.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.984.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.367.49'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('E20').Value = '  -2.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.99'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.59'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('E24').Value = '  -0.36%  '
$ws.Range('E26').Value = '  +1.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E32').Value = '  -4.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('E35').Value = '  -4.48%  '
$ws.Range('E36').Value = '  -1.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.89'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('E38').Value = '  -3.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.22%  '
$ws.Range('E40').Value = '  +0.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.685.10'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('E44').Value = '  -2.85%  '
$ws.Range('E45').Value = '  -2.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '335.87'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.02%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.32'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('E50').Value = '  +2.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.04%  '
